# Visualization now places dots logically in higher-dimensional problems.
#
# - GA-innstillinger (sheet 1) gains three explanatory notes about how the
#   bitstring->coordinate mapping and distance measure work, and a note
#   about the new peak-detection epsilon, and becomes the active sheet.
# - Eksperiment 1 - Klyngingsalgori (sheet 2) gets its two existing notes
#   pushed down two rows and a new results table (Function / Dimensionality
#   / Best / Worst / Mean / STD header row plus a couple of data rows)
#   added beneath them.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1: GA-innstillinger (part 1) -----------------------------------
$ws1.Range("A23").Value = "Mapping bitstrings to numbers by parsing them as binary numbers, then normalizing to the search space."
$ws1.Range("A24").Value = "Euclidean distance measure on the normalized axes."

# --- Sheet 2: Eksperiment 1 - Klyngingsalgori -----------------------------
# Push the two existing rows down from 1,2 to 3,4 by inserting two rows at
# the top (keeps their formatting intact).
$ws2.Rows.Item(1).Insert()
$ws2.Rows.Item(1).Insert()

# New results table header (row 6, bold) ...
$ws2.Range("A6").Value = "Function"
$ws2.Range("B6").Value = "Dimensionality"
$ws2.Range("C6").Value = "Best"
$ws2.Range("D6").Value = "Worst"
$ws2.Range("E6").Value = "Mean"
$ws2.Range("F6").Value = "STD"
$ws2.Range("A6:F6").Font.Bold = $true

# ... and a few data rows beneath it.
$ws2.Range("A7").Value = 1
$ws2.Range("B7").Value = 5
$ws2.Range("A8").Value = 1
$ws2.Range("A9").Value = 1

# --- Sheet 1: GA-innstillinger (part 2) -----------------------------------
$ws1.Range("A26").Value = "Peak detection epsilon is 0.5*D"

# --- Selection / active sheet --------------------------------------------
# GA-innstillinger becomes the active (tab-selected) sheet with A27 selected,
# while the experiment sheet's selection moves to C7.
$ws2.Range("C7").Select() | Out-Null
$ws1.Select() | Out-Null
$ws1.Range("A27").Select() | Out-Null
